# Atualizado por script em 21-11-2023 14:30
#
# The spreadsheet held several rows whose (home/away/odds/url) details
# had been assigned to the wrong fixture row within their same-date
# block. This script fixes the mix-ups by re-distributing the F:V data
# among the affected rows, and appends the newly scraped fixture
# (Orleta Radzyn x Chelmianka Chelm) as row 134.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowValues($row) {
    $result = @{}
    foreach ($c in $cols) {
        $result[$c] = $ws.Range($c + $row).Value()
    }
    return $result
}

function Set-RowValues($row, $data) {
    foreach ($c in $cols) {
        $ws.Range($c + $row).Value = $data[$c]
    }
}

# ---------------------------------------------------------------------
# 1) Snapshot the current (pre-fix) F:V contents of every affected row
#    before any writes happen, so later writes can't clobber a value
#    that is still needed as a source for another row.
# ---------------------------------------------------------------------
$affectedRows = @(66,67,79,80,81,82,89,90,91,113,114,115,116,120,121,122)
$snapshot = @{}
foreach ($r in $affectedRows) {
    $snapshot[$r] = Get-RowValues $r
}

# ---------------------------------------------------------------------
# 2) new[row] <- old[sourceRow] : re-assign the F:V block of data.
#    Columns A:E (index/country/tournament/season/date) are untouched.
# ---------------------------------------------------------------------
$mapping = @{
    66  = 67;  67  = 66;
    79  = 82;  80  = 79;  81  = 80;  82  = 81;
    89  = 91;  90  = 89;  91  = 90;
    113 = 114; 114 = 113;
    115 = 116; 116 = 115;
    120 = 121; 121 = 122; 122 = 120;
}

foreach ($newRow in $mapping.Keys) {
    $sourceRow = $mapping[$newRow]
    Set-RowValues $newRow $snapshot[$sourceRow]
}

# ---------------------------------------------------------------------
# 3) Append the new fixture as row 134, copying row 133's formatting
#    (bold/bordered index style in column A, date format in column E)
#    before filling in the values.
# ---------------------------------------------------------------------
$ws.Range("A133:V133").Copy()
$ws.Range("A134:V134").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A134").Value = 133
$ws.Range("B134").Value = "poland"
$ws.Range("C134").Value = "iii-liga-group-iv"
$ws.Range("D134").Value = "2023-2024"
$ws.Range("E134").Value = 45251.54166666666
$ws.Range("F134").Value = "Orleta Radzyn"
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = "Chelmianka Chelm"
$ws.Range("I134").Value = 2
$ws.Range("J134").Value = 3.42
$ws.Range("K134").Value = "18/11/2023 01:13"
$ws.Range("L134").Value = 3.13
$ws.Range("M134").Value = "21/11/2023 12:59"
$ws.Range("N134").Value = 3.65
$ws.Range("O134").Value = "18/11/2023 01:13"
$ws.Range("P134").Value = 3.58
$ws.Range("Q134").Value = "21/11/2023 12:59"
$ws.Range("R134").Value = 1.8
$ws.Range("S134").Value = "18/11/2023 01:13"
$ws.Range("T134").Value = 1.97
$ws.Range("U134").Value = "21/11/2023 12:59"
$ws.Range("V134").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iv/orleta-radzyn-chelmianka-chelm/dM8HHVcs/"
